# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) timestamps on the
# per-language handback sheets to reflect the newly generated report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 12:26:09"
$wsZhCn.Range("E5").Value = "2016-03-22 12:26:09"
$wsZhCn.Range("H2").Value = "2016-03-22 12:26:37"
$wsZhCn.Range("H5").Value = "2016-03-22 12:26:37"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 12:26:13"
$wsDeDe.Range("E5").Value = "2016-03-22 12:26:13"
$wsDeDe.Range("H2").Value = "2016-03-22 12:26:43"
$wsDeDe.Range("H5").Value = "2016-03-22 12:26:43"
